# Updates cryptos list price/volume columns (and a few re-ranked rows whose
# Coin/Link/Price/Volume values moved to a different row) to match the
# latest scrape. Numeric-looking Price values are written with a leading
# apostrophe so Excel stores them as text (quote-prefix) instead of
# auto-converting to a number, matching the sheet's original text-cell
# formatting for the Price column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.506.07'
$ws.Range('E2').Value = '  +1.15%  '

$ws.Range('D3').Value = '1.795.68'
$ws.Range('E3').Value = '  +1.07%  '

$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '''227.31'
$ws.Range('E5').Value = '  +1.02%  '

$ws.Range('D6').Value = '''0.556'
$ws.Range('E6').Value = '  +2.16%  '

$ws.Range('E7').Value = '  -0.29%  '

$ws.Range('D8').Value = '''32.92'
$ws.Range('E8').Value = '  +4.22%  '

$ws.Range('D9').Value = '''0.296'
$ws.Range('E9').Value = '  +2.02%  '

$ws.Range('D10').Value = '''0.0694'
$ws.Range('E10').Value = '  +1.68%  '

$ws.Range('D11').Value = '''0.0947'
$ws.Range('E11').Value = '  +0.08%  '

$ws.Range('D12').Value = '2.056.84'
$ws.Range('E12').Value = '  +1.21%  '

$ws.Range('D13').Value = '''11.16'
$ws.Range('E13').Value = '  +2.93%  '

$ws.Range('D14').Value = '1.786.95'
$ws.Range('E14').Value = '  -0.39%  '

$ws.Range('D15').Value = '''0.640'
$ws.Range('E15').Value = '  +3.57%  '

$ws.Range('D16').Value = '34.521.16'
$ws.Range('E16').Value = '  +1.25%  '

$ws.Range('E17').Value = '  +3.33%  '

$ws.Range('D18').Value = '''69.08'
$ws.Range('E18').Value = '  +2.26%  '

$ws.Range('D19').Value = '0.0₃0805'
$ws.Range('E19').Value = '  +2.71%  '

$ws.Range('D20').Value = '''246.22'
$ws.Range('E20').Value = '  +0.81%  '

$ws.Range('D21').Value = '''11.37'
$ws.Range('E21').Value = '  +3.59%  '

$ws.Range('E22').Value = '  -0.40%  '

$ws.Range('D23').Value = '''4.17'
$ws.Range('E23').Value = '  +2.28%  '

$ws.Range('D24').Value = '''170.70'
$ws.Range('E24').Value = '  +6.00%  '

$ws.Range('D25').Value = '''2.06'
$ws.Range('E25').Value = '  +1.45%  '

$ws.Range('D26').Value = '''7.36'
$ws.Range('E26').Value = '  +3.59%  '

$ws.Range('D27').Value = '''16.67'
$ws.Range('E27').Value = '  +3.04%  '

$ws.Range('D28').Value = '''0.116'
$ws.Range('E28').Value = '  +2.25%  '

$ws.Range('E29').Value = '  -0.40%  '

$ws.Range('D30').Value = '''4.02'
$ws.Range('E30').Value = '  +8.74%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.0526'
$ws.Range('E31').Value = '  +2.16%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''3.82'
$ws.Range('E32').Value = '  +3.37%  '

$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '''1.24'
$ws.Range('E33').Value = '  +0.99%  '

$ws.Range('D34').Value = '''1.83'
$ws.Range('E34').Value = '  +2.85%  '

$ws.Range('D35').Value = '1.419.59'
$ws.Range('E35').Value = '  -1.14%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''0.684'
$ws.Range('E36').Value = '  +5.10%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''2.56'
$ws.Range('E37').Value = '  +6.61%  '

$ws.Range('E38').Value = '  +3.58%  '

$ws.Range('E39').Value = '  +0.31%  '

$ws.Range('D40').Value = '''84.45'
$ws.Range('E40').Value = '  +5.72%  '

$ws.Range('D41').Value = '''0.948'
$ws.Range('E41').Value = '  +3.76%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '''2.78'
$ws.Range('E42').Value = '  +2.07%  '

$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Value = '''2.39'
$ws.Range('E43').Value = '  +0.57%  '

$ws.Range('D44').Value = '''14.00'
$ws.Range('E44').Value = '  +4.62%  '

$ws.Range('D45').Value = '''0.0528'
$ws.Range('E45').Value = '  +2.63%  '

$ws.Range('D46').Value = '''1.10'
$ws.Range('E46').Value = '  +2.13%  '

$ws.Range('D47').Value = '''6.13'
$ws.Range('E47').Value = '  +1.71%  '

$ws.Range('D48').Value = '1.957.07'
$ws.Range('E48').Value = '  +1.14%  '

$ws.Range('D49').Value = '''105.45'
$ws.Range('E49').Value = '  +1.54%  '

$ws.Range('E50').Value = '  -0.32%  '

$ws.Range('E51').Value = '  -1.67%  '
